# Increment the "想去人数" (interested-count) column F for a handful of
# events on the "展览" and "全部类型" sheets, matching the upstream data
# refresh captured in the diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row => new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(4, 6).Value = 553
$ws1.Cells.Item(16, 6).Value = 953
$ws1.Cells.Item(18, 6).Value = 1612
$ws1.Cells.Item(19, 6).Value = 69
$ws1.Cells.Item(28, 6).Value = 1040
$ws1.Cells.Item(32, 6).Value = 170
$ws1.Cells.Item(33, 6).Value = 1554
$ws1.Cells.Item(34, 6).Value = 2092
$ws1.Cells.Item(35, 6).Value = 990
$ws1.Cells.Item(36, 6).Value = 30
$ws1.Cells.Item(37, 6).Value = 232
$ws1.Cells.Item(42, 6).Value = 361

# Sheet "全部类型": row => new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(7, 6).Value = 553
$ws4.Cells.Item(18, 6).Value = 953
$ws4.Cells.Item(20, 6).Value = 1612
$ws4.Cells.Item(22, 6).Value = 69
$ws4.Cells.Item(34, 6).Value = 1040
$ws4.Cells.Item(36, 6).Value = 1554
$ws4.Cells.Item(37, 6).Value = 2092
$ws4.Cells.Item(39, 6).Value = 990
$ws4.Cells.Item(40, 6).Value = 30
$ws4.Cells.Item(42, 6).Value = 232
$ws4.Cells.Item(46, 6).Value = 361
